$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("missed")

$values = @(
    "40-2-2",
    "46-1-2",
    "30-4-3",
    "45-5-2",
    "66-3-1",
    "46-3-2",
    "47-1-1",
    "45-5-1",
    "37-5-1",
    "50-1-2",
    "49-2-1",
    "38-2-1",
    "50-3-2",
    "49-3-2",
    "46-3-1",
    "37-1-1",
    "41-3-1",
    "40-4-2",
    "56-2-2",
    "49-1-1",
    "51-2-2",
    "50-4-1",
    "48-4-1",
    "49-2-2",
    "52-3-1",
    "52-2-1",
    "50-3-1",
    "42-3-2",
    "38-4-1",
    "46-4-1",
    "52-2-2",
    "40-5-1",
    "49-1-2",
    "47-4-2",
    "38-1-1",
    "50-2-1",
    "51-1-1"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
